$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calibration params")

$ws.Range("B4").Value = 2.0487476081399998
$ws.Range("C4").Value = 2.3317645809899998
$ws.Range("D4").Value = 1.9793818972999999
$ws.Range("E4").Value = 2.1467547205400002
$ws.Range("F4").Value = 1.8948320172999999
$ws.Range("G4").Value = 2.0401211340800001
$ws.Range("H4").Value = 2.1565290911099999
$ws.Range("I4").Value = 2.1899351025099998

$ws.Range("B4:I4").Select()
